$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$school = "โรงเรียนโคกเพชรวิทยาคาร"
$detail = "วัสดุการศึกษา"
$note = "Created on 14-01-2024"

$rows = @(
    @{ Row = 26; A = 22; D = 1 },
    @{ Row = 27; A = 23; D = 5 },
    @{ Row = 28; A = 24; D = 1 }
)

foreach ($r in $rows) {
    $ws.Range("A" + $r.Row).Value = $r.A
    $ws.Range("B" + $r.Row).Value = $school
    $ws.Range("C" + $r.Row).Value = $detail
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("E" + $r.Row).Value = $note
}
